$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D7").Value = -7.232099999999998
$ws.Range("C8").Value = -12.28499999999999
$ws.Range("C10").Value = -13.12289999999999
$ws.Range("C12").Value = -10.8028
$ws.Range("D15").Value = -7.963799999999997
$ws.Range("C18").Value = -14.01479999999999
$ws.Range("D18").Value = -8.530599999999993
$ws.Range("D20").Value = -7.912899999999994
$ws.Range("D29").Value = -7.315700000000001
$ws.Range("D30").Value = -7.288400000000001
$ws.Range("D31").Value = -8.537499999999996
$ws.Range("C37").Value = -12.80400000000001
$ws.Range("D40").Value = -8.157799999999991
$ws.Range("D50").Value = -8.146499999999998
$ws.Range("C55").Value = -14.0446
$ws.Range("C68").Value = -10.9615
$ws.Range("D68").Value = -7.165399999999994
$ws.Range("D76").Value = -7.2878
$ws.Range("C77").Value = -12.4132
$ws.Range("C78").Value = -12.51970000000001
$ws.Range("C81").Value = -13.0803
$ws.Range("C82").Value = -12.322
$ws.Range("D87").Value = -7.938299999999995
$ws.Range("D88").Value = -7.462399999999998
$ws.Range("D96").Value = -7.420500000000003
$ws.Range("D98").Value = -8.296600000000005
$ws.Range("D101").Value = -7.758699999999999
$ws.Range("D102").Value = -7.727399999999998
